$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.326.04'
$ws.Range("E2").Value = '  +5.53%  '
$ws.Range("D3").Value = '3.391.22'
$ws.Range("E3").Value = '  +6.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.20'
$ws.Range("E5").Value = '  +7.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.72'
$ws.Range("E6").Value = '  +6.57%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.392.18'
$ws.Range("E8").Value = '  +6.14%  '
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.48'
$ws.Range("E10").Value = '  +2.12%  '
$ws.Range("E11").Value = '  +6.84%  '
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("D13").Value = '3.972.00'
$ws.Range("E13").Value = '  +6.10%  '
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("E15").Value = '  +6.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.01'
$ws.Range("E16").Value = '  +4.71%  '
$ws.Range("D17").Value = '63.427.03'
$ws.Range("E17").Value = '  +5.65%  '
$ws.Range("D18").Value = '3.387.99'
$ws.Range("E18").Value = '  +5.94%  '
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.98'
$ws.Range("E20").Value = '  +5.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.42'
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.72'
$ws.Range("E22").Value = '  +5.31%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.535'
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.79'
$ws.Range("E25").Value = '  +2.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.52'
$ws.Range("E26").Value = '  +11.22%  '
$ws.Range("E27").Value = '  +6.21%  '
$ws.Range("E28").Value = '  +18.41%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +7.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.46'
$ws.Range("E31").Value = '  +5.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.12'
$ws.Range("E32").Value = '  +2.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.32'
$ws.Range("E33").Value = '  +10.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.59'
$ws.Range("E34").Value = '  +5.64%  '
$ws.Range("E35").Value = '  +2.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.49'
$ws.Range("E36").Value = '  +9.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '158.57'
$ws.Range("E37").Value = '  +1.41%  '
$ws.Range("E38").Value = '  +11.96%  '
$ws.Range("E39").Value = '  +4.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0751'
$ws.Range("E40").Value = '  +6.48%  '
$ws.Range("D41").Value = '2.879.58'
$ws.Range("E41").Value = '  +3.18%  '
$ws.Range("E42").Value = '  +4.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.762'
$ws.Range("E43").Value = '  +5.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.20'
$ws.Range("E44").Value = '  +3.91%  '
$ws.Range("E45").Value = '  +1.18%  '
$ws.Range("E46").Value = '  +7.97%  '
$ws.Range("D47").Value = '3.433.83'
$ws.Range("E47").Value = '  +6.20%  '
$ws.Range("E48").Value = '  +6.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '298.66'
$ws.Range("E49").Value = '  +12.81%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.33'
$ws.Range("E50").Value = '  +3.11%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.103'
$ws.Range("E51").Value = '  -0.86%  '
